$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names (Column B) for the 24 data rows, in order.
$names = @("John","Sam","Ricardo","Jose","Joanna","Juan","Luis","Tomas","Matthew","Christ","Leonardo","Jennifer","Angelina","Natasha","Silva","Jonathan","Andres","Andrea","Helena","Ariadna","Rodrigo","Antonio","Philip","Bastian")

# Ages (Column C) for the 24 data rows, in order.
$ages = @(25,19,34,35,25,26,21,24,28,31,35,30,31,34,21,31,30,29,32,18,29,23,28,22)

# Resize the table to include the new header row plus all 24 data rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C25"))

# Write the header row.
$ws.Cells.Item(1, 1).Value = "Id"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "Age"

# Write each data row (shifted down one row to make room for the header).
for ($i = 0; $i -lt 24; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $ages[$i]
}
